$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Update selection to B6
$ws.Range("B6").Select()

# Compilation success row: B5 -> "no", C5 -> note
$ws.Range("B5").Value = "no"
$ws.Range("C5").Value = "Wrong page object model name"

# Runtime without error row: clear B6 value
$ws.Range("B6").Value = $null

# Assertion validity row: clear B7 and C7 values
$ws.Range("B7").Value = $null
$ws.Range("C7").Value = $null

# Code BLEU row: update score value and note text
$ws.Range("B12").Value = 0.3019884666684997
$ws.Range("C12").Value = "{'codebleu': 0.30198846666849966, 'ngram_match_score': 0.0922259976907841, 'weighted_ngram_match_score': 0.10259829515626676, 'syntax_match_score': 0.6435643564356436, 'dataflow_match_score': 0.3695652173913043}"

$wb.Save()
